$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" contain duplicate event listings that both
# need their "想去人数" (want-to-go count) column (F) refreshed.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 461
    $ws.Range("F3").Value = 3286
    $ws.Range("F5").Value = 655
}
